$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week1")

$ws.Range("C9").Value = 0.010416666666666666
$ws.Range("C10").Value = 0.006944444444444444

$ws.Range("C11").Select()
